# Auto-generated: update Thbs1-Itga6 NATMI TPM values per commit "update scripts wuth new tpm"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.727484333333333
$ws.Range("H2").Value = 8.182453000000001
$ws.Range("I2").Value = 0.03096049453772388
$ws.Range("J2").Value = 0.03096049453772388
$ws.Range("M2").Value = 133.3951123333333
$ws.Range("N2").Value = 400.185337
$ws.Range("O2").Value = 0.8984588679103155
$ws.Range("P2").Value = 0.8984588679103156
$ws.Range("Q2").Value = 363.8330790324068
$ws.Range("R2").Value = 3274.497711291661
$ws.Range("S2").Value = 0.02781673087230691
$ws.Range("T2").Value = 0.02781673087230691

# Row 3
$ws.Range("G3").Value = 2.727484333333333
$ws.Range("H3").Value = 8.182453000000001
$ws.Range("I3").Value = 0.03096049453772388
$ws.Range("J3").Value = 0.03096049453772388
$ws.Range("M3").Value = 2.340788333333334
$ws.Range("N3").Value = 7.022365000000001
$ws.Range("O3").Value = 0.01576596023045448
$ws.Range("P3").Value = 0.01576596023045448
$ws.Range("Q3").Value = 6.384463506816112
$ws.Range("R3").Value = 57.46017156134501
$ws.Range("S3").Value = 0.0004881219255969578
$ws.Range("T3").Value = 0.000488121925596958

# Row 4
$ws.Range("G4").Value = 2.727484333333333
$ws.Range("H4").Value = 8.182453000000001
$ws.Range("I4").Value = 0.03096049453772388
$ws.Range("J4").Value = 0.03096049453772388
$ws.Range("M4").Value = 12.735128
$ws.Range("N4").Value = 38.205384
$ws.Range("O4").Value = 0.08577517185923002
$ws.Range("P4").Value = 0.08577517185923003
$ws.Range("Q4").Value = 34.73486210299467
$ws.Range("R4").Value = 312.613758926952
$ws.Range("S4").Value = 0.002655641739820018
$ws.Range("T4").Value = 0.002655641739820019

# Row 5
$ws.Range("I5").Value = 0.5986009007423507
$ws.Range("J5").Value = 0.5986009007423507
$ws.Range("M5").Value = 133.3951123333333
$ws.Range("N5").Value = 400.185337
$ws.Range("O5").Value = 0.8984588679103155
$ws.Range("P5").Value = 0.8984588679103156
$ws.Range("Q5").Value = 7034.474483710005
$ws.Range("R5").Value = 63310.27035339003
$ws.Range("S5").Value = 0.5378182876110675
$ws.Range("T5").Value = 0.5378182876110676

# Row 6
$ws.Range("I6").Value = 0.5986009007423507
$ws.Range("J6").Value = 0.5986009007423507
$ws.Range("M6").Value = 2.340788333333334
$ws.Range("N6").Value = 7.022365000000001
$ws.Range("O6").Value = 0.01576596023045448
$ws.Range("P6").Value = 0.01576596023045448
$ws.Range("S6").Value = 0.00943751799501813
$ws.Range("T6").Value = 0.009437517995018132

# Row 7
$ws.Range("I7").Value = 0.5986009007423507
$ws.Range("J7").Value = 0.5986009007423507
$ws.Range("M7").Value = 12.735128
$ws.Range("N7").Value = 38.205384
$ws.Range("O7").Value = 0.08577517185923002
$ws.Range("P7").Value = 0.08577517185923003
$ws.Range("Q7").Value = 671.5758275979574
$ws.Range("R7").Value = 6044.182448381616
$ws.Range("S7").Value = 0.05134509513626502
$ws.Range("T7").Value = 0.05134509513626503

# Row 8
$ws.Range("G8").Value = 32.63402300000001
$ws.Range("H8").Value = 97.90206900000001
$ws.Range("I8").Value = 0.3704386047199253
$ws.Range("J8").Value = 0.3704386047199253
$ws.Range("M8").Value = 133.3951123333333
$ws.Range("N8").Value = 400.185337
$ws.Range("O8").Value = 0.8984588679103155
$ws.Range("P8").Value = 0.8984588679103156
$ws.Range("Q8").Value = 4353.219163973585
$ws.Range("R8").Value = 39178.97247576226
$ws.Range("S8").Value = 0.3328238494269409
$ws.Range("T8").Value = 0.332823849426941

# Row 9
$ws.Range("G9").Value = 32.63402300000001
$ws.Range("H9").Value = 97.90206900000001
$ws.Range("I9").Value = 0.3704386047199253
$ws.Range("J9").Value = 0.3704386047199253
$ws.Range("M9").Value = 2.340788333333334
$ws.Range("N9").Value = 7.022365000000001
$ws.Range("O9").Value = 0.01576596023045448
$ws.Range("P9").Value = 0.01576596023045448
$ws.Range("Q9").Value = 76.38934030813169
$ws.Range("R9").Value = 687.5040627731852
$ws.Range("S9").Value = 0.005840320309839388
$ws.Range("T9").Value = 0.00584032030983939

# Row 10
$ws.Range("G10").Value = 32.63402300000001
$ws.Range("H10").Value = 97.90206900000001
$ws.Range("I10").Value = 0.3704386047199253
$ws.Range("J10").Value = 0.3704386047199253
$ws.Range("M10").Value = 12.735128
$ws.Range("N10").Value = 38.205384
$ws.Range("O10").Value = 0.08577517185923002
$ws.Range("P10").Value = 0.08577517185923003
$ws.Range("Q10").Value = 415.5984600599441
$ws.Range("R10").Value = 3740.386140539497
$ws.Range("S10").Value = 0.03177443498314497
$ws.Range("T10").Value = 0.03177443498314497

